# Apply referee stats update for 2025-11-22 publish run.
$wb = $excel.ActiveWorkbook

$timestamp = "2025-11-22 03:03:39"

# --- Sheet "Главные" (head referees) ---
$wsMain = $wb.Worksheets.Item("Главные")

# Row 9 - Gashilov Viktor
$wsMain.Range("C9").Value = 27
$wsMain.Range("D9").Value = 420
$wsMain.Range("E9").Value = 225
$wsMain.Range("F9").Value = 195
$wsMain.Range("G9").Value = 15.56
$wsMain.Range("H9").Value = 8.33
$wsMain.Range("I9").Value = 7.22
$wsMain.Range("J9").Value = 110
$wsMain.Range("K9").Value = 95
$wsMain.Range("W9").Value = 24

# Row 14 - Lavrentev Anton
$wsMain.Range("C14").Value = 19
$wsMain.Range("D14").Value = 244
$wsMain.Range("E14").Value = 127
$wsMain.Range("F14").Value = 117
$wsMain.Range("G14").Value = 12.84
$wsMain.Range("H14").Value = 6.68
$wsMain.Range("I14").Value = 6.16
$wsMain.Range("J14").Value = 61
$wsMain.Range("K14").Value = 51
$wsMain.Range("W14").Value = 12

# Row 18 - Naumov Denis
$wsMain.Range("C18").Value = 26
$wsMain.Range("D18").Value = 391
$wsMain.Range("E18").Value = 182
$wsMain.Range("F18").Value = 209
$wsMain.Range("G18").Value = 15.04
$wsMain.Range("H18").Value = 7
$wsMain.Range("I18").Value = 8.039999999999999
$wsMain.Range("J18").Value = 81
$wsMain.Range("K18").Value = 97

# Row 21 - Romasko Evgeniy
$wsMain.Range("C21").Value = 23
$wsMain.Range("D21").Value = 320
$wsMain.Range("E21").Value = 144
$wsMain.Range("F21").Value = 176
$wsMain.Range("G21").Value = 13.91
$wsMain.Range("H21").Value = 6.26
$wsMain.Range("I21").Value = 7.65
$wsMain.Range("J21").Value = 62
$wsMain.Range("K21").Value = 73

# Refresh as_of_utc stamp for all data rows (2-26)
for ($r = 2; $r -le 26; $r++) {
    $wsMain.Range("AA$r").Value = $timestamp
}

# --- Sheet "Линейные" (line referees) ---
$wsLine = $wb.Worksheets.Item("Линейные")

# Row 17 - Litvinov Aleksandr
$wsLine.Range("C17").Value = 14
$wsLine.Range("D17").Value = 264
$wsLine.Range("E17").Value = 153
$wsLine.Range("F17").Value = 111
$wsLine.Range("G17").Value = 18.86
$wsLine.Range("H17").Value = 10.93
$wsLine.Range("I17").Value = 7.93
$wsLine.Range("J17").Value = 59
$wsLine.Range("K17").Value = 48
$wsLine.Range("W17").Value = 6

# Row 22 - Sedov Egor
$wsLine.Range("C22").Value = 20
$wsLine.Range("D22").Value = 387
$wsLine.Range("E22").Value = 199
$wsLine.Range("F22").Value = 188
$wsLine.Range("G22").Value = 19.35
$wsLine.Range("H22").Value = 9.949999999999999
$wsLine.Range("I22").Value = 9.4
$wsLine.Range("J22").Value = 82
$wsLine.Range("K22").Value = 84
$wsLine.Range("W22").Value = 24

# Refresh as_of_utc stamp for all data rows (2-26)
for ($r = 2; $r -le 26; $r++) {
    $wsLine.Range("AA$r").Value = $timestamp
}
